$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new header columns (P1, Q1), continuing the 0..13 sequence in row 1,
# and copy the header formatting (bold, centered, bordered) from the existing
# header cell O1 so the new cells pick up the same style.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null

# Update existing columns I, K, M, O for the data rows (2-25), and populate the
# two newly added columns P and Q.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P = 2
    $ws.Cells.Item($r, 17).Value = 2  # Q = 2
}
